$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I17").Value = 43734
